# fix: revert admin dev default; seed customers only when table empty;
# autosave on customer select when hours/day present
#
# Data values in the two sample sheets are being reset back to "empty
# defaults": the seeded sample client names / employee id are swapped
# out for fresh placeholder data, the PTO "seeded sample" rows become
# plain "Regular" rows, the Notes column is cleared, and every Rate /
# Total cell drops back to 0 (no more hard-coded $100 / $800 demo
# numbers).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "Weekly Timesheet"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Weekly Timesheet")

# New client names for the 5 seeded rows (rows 2-6)
$clients = @("Smithers", "Bottomley", "Behrens", "Goodrich", "Campbell")
for ($i = 0; $i -lt $clients.Length; $i++) {
    $row = 2 + $i
    $ws1.Cells.Item($row, 2).Value = $clients[$i]   # B: Client
    $ws1.Cells.Item($row, 4).Value = "Regular"       # D: Type (was PTO/Regular)
    $ws1.Cells.Item($row, 5).Value = 0               # E: Rate
    $ws1.Cells.Item($row, 6).Value = 0               # F: Total
}

# Subtotal / admin-subtotal / grand-total rows no longer carry the
# hard-coded $4000 demo figure.
$ws1.Cells.Item(8, 6).Value = 0    # F8  SUBTOTAL total
$ws1.Cells.Item(12, 6).Value = 0   # F12 ADMIN SUBTOTAL total
$ws1.Cells.Item(13, 6).Value = 0   # F13 GRAND TOTAL total

# ---------------------------------------------------------------
# Sheet 2: "Jason Schema"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Jason Schema")

# Employee ID reverted to the non-seeded default id
for ($row = 2; $row -le 6; $row++) {
    $ws2.Cells.Item($row, 2).Value = "emp_lf0u97k0"   # B: Employee ID
}

for ($i = 0; $i -lt $clients.Length; $i++) {
    $row = 2 + $i
    $ws2.Cells.Item($row, 4).Value = $clients[$i]   # D: Client
    $ws2.Cells.Item($row, 6).Value = 0               # F: Rate
    $ws2.Cells.Item($row, 7).Value = 0               # G: Total
    $ws2.Cells.Item($row, 8).Value = "Regular"       # H: Type (was PTO/Regular)
    $ws2.Cells.Item($row, 9).Value = ""              # I: Notes (seeded-sample text cleared)
}
